$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data added beneath the existing CR entries (row 6), documenting
# CR1.1 - "Diplaying the list of Proforma(Letters) without search, while creating".
# Cells are written in the same order the strings were originally added to the
# workbook's shared-string table so the resulting index order matches.

$ws.Range("A6").Value = "March 5,2017"

$sql = "INSERT INTO FunctionIdentifier (fidName, projectPackage, codePackage, description, fidEnabled) values ('MOLSACommunicationDA.listTemplateByTypeAndParticipantAutoDisplay', 'curam', 'molsa.core.facade', 'curam.molsa.core.facade.impl.MOLSACommunicationDA.listTemplateByTypeAndParticipantAutoDisplay', 'Y');`nINSERT INTO SECURITYIDENTIFIER (DESCRIPTION, LASTWRITTEN, SIDNAME, SIDTYPE, VERSIONNO) VALUES (null, '2017-03-05 03:42:25', 'MOLSACommunicationDA.listTemplateByTypeAndParticipantAutoDisplay', 'FUNCTION', 0);`nINSERT INTO SECURITYFIDSID(SIDNAME, FIDNAME) VALUES ('MOLSACommunicationDA.listTemplateByTypeAndParticipantAutoDisplay','MOLSACommunicationDA.listTemplateByTypeAndParticipantAutoDisplay');`nINSERT INTO SECURITYGROUPSID (GROUPNAME, LASTWRITTEN, SIDNAME) VALUES ('SUPERGROUP', null, 'MOLSACommunicationDA.listTemplateByTypeAndParticipantAutoDisplay');"
$ws.Range("G6").Value = $sql
$ws.Range("G6").WrapText = $true

$ws.Range("E6").Value = "Joseph"

$ws.Range("C6").Value = "/EJBServer/components/MOLSA/model/Packages/Reference Model/Facade/Communication.efx`n"

$ws.Range("D6").Value = "CR1.1 Diplaying the list of Proforma(Letters) without search, while creating"
$ws.Range("D6").WrapText = $true

$ws.Range("F6").Value = "Yes"

# Row 6 grew tall to fit the wrapped SQL/comment text.
$ws.Rows.Item(6).RowHeight = 150

# Move the active selection the way the author's workbook view ended up.
$ws.Range("D10").Select() | Out-Null
